$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8643032312393188
$ws.Range("B1").Value = 2.626749515533447
$ws.Range("C1").Value = 3.354634523391724
$ws.Range("D1").Value = 1.886842250823975
$ws.Range("E1").Value = 1.444468140602112
